# Updated cryptos list on Fri Sep 15 03:48:44 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns with the latest
# scraped coinranking.com figures. A leading apostrophe is used on the
# Price cells whose new text happens to look like a plain number (e.g.
# "212.99") so Excel keeps storing them as literal text, matching the
# existing text values elsewhere in the column (e.g. "26.644.38").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.644.38"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "1.636.65"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'212.99"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  +2.52%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("D11").Value = "'0.0840"
$ws.Range("E11").Value = "  +3.06%  "
$ws.Range("D12").Value = "1.864.60"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").Value = "1.630.38"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("E14").Value = "  +2.15%  "
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").Value = "26.661.88"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").Value = "'63.37"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "0.0₃0745"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").Value = "'218.82"
$ws.Range("E19").Value = "  +7.76%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'4.30"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "'9.49"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("E23").Value = "  +2.64%  "
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").Value = "'148.69"
$ws.Range("E25").Value = "  +3.70%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  +4.34%  "
$ws.Range("D29").Value = "'15.42"
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("E30").Value = "  -2.87%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  +4.07%  "
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("E35").Value = "  -1.57%  "
$ws.Range("D36").Value = "1.195.40"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").Value = "'0.0174"
$ws.Range("E37").Value = "  +5.75%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "'5.41"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("D43").Value = "'0.793"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("D44").Value = "1.774.21"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("D45").Value = "'92.15"
$ws.Range("E45").Value = "  -1.56%  "
$ws.Range("D47").Value = "'54.78"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").Value = "'0.0512"
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("D49").Value = "'7.65"
$ws.Range("E49").Value = "  +5.10%  "
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("E51").Value = "  +0.04%  "
